$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 31749.75
$ws.Range("I21").Value = 31749.75
$ws.Range("K21").Value = 31749.75
$ws.Range("M21").Value = -31281.75

$ws.Range("H23").Value = 31749.75
$ws.Range("I23").Value = 31749.75
$ws.Range("K23").Value = 31749.75
$ws.Range("M23").Value = -31515.75

$ws.Range("H28").Value = 5707.273
$ws.Range("I28").Value = 2784.75
$ws.Range("K28").Value = 2784.75
$ws.Range("M28").Value = -2299.75

$ws.Range("H51").Value = 7617.087
$ws.Range("I51").Value = 20231
$ws.Range("K51").Value = 20231
$ws.Range("M51").Value = -19747

$ws.Range("H53").Value = 66667670
$ws.Range("J53").Value = 958.7
$ws.Range("L53").Value = 958.7
$ws.Range("N53").Value = -2232.7

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4431.0576
$ws.Range("I32").Value = 4557.1304
$ws.Range("K32").Value = 4557.1304
$ws.Range("M32").Value = -4270.1304

$ws.Range("H43").Value = 45075.4
$ws.Range("J43").Value = 48995.668
$ws.Range("L43").Value = 48995.668
$ws.Range("N43").Value = -49621.668

$ws.Range("H61").Value = 7243708.5
$ws.Range("I61").Value = 7145269.5
$ws.Range("K61").Value = 7145269.5
$ws.Range("M61").Value = -7145057.5

$ws.Range("H74").Value = 695825.4399999999
$ws.Range("I74").Value = 758506.1
$ws.Range("K74").Value = 758506.1
$ws.Range("M74").Value = -757632.1

$ws.Range("H77").Value = 695825.4399999999
$ws.Range("I77").Value = 758506.1
$ws.Range("K77").Value = 3792530.5
$ws.Range("M77").Value = -3788162.5

$ws.Range("H110").Value = 7519.4165
$ws.Range("I110").Value = 7023.3
$ws.Range("K110").Value = 7023.3
$ws.Range("M110").Value = -4978.3

$ws.Range("H132").Value = 2003009.6
$ws.Range("I132").Value = 2793.111
$ws.Range("K132").Value = 8379.332999999999
$ws.Range("M132").Value = -5849.332999999999

$ws.Range("H136").Value = 7243708.5
$ws.Range("I136").Value = 7145269.5
$ws.Range("K136").Value = 21435808.5
$ws.Range("M136").Value = -21433258.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4173.55
$ws.Range("I107").Value = 3815.111
$ws.Range("K107").Value = 3815.111
$ws.Range("M107").Value = -1895.111

$ws.Range("H137").Value = 98784.5
$ws.Range("J137").Value = 98784.5
$ws.Range("L137").Value = 98784.5
$ws.Range("N137").Value = -108984.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 7149277
$ws.Range("I16").Value = 10005296
$ws.Range("K16").Value = 10005296
$ws.Range("M16").Value = -10005009

$ws.Range("H31").Value = 40406824
$ws.Range("I31").Value = 50002736
$ws.Range("K31").Value = 50002736
$ws.Range("M31").Value = -50002441

$ws.Range("H34").Value = 40406824
$ws.Range("I34").Value = 50002736
$ws.Range("K34").Value = 50002736
$ws.Range("M34").Value = -50002534

$ws.Range("H113").Value = 7149277
$ws.Range("I113").Value = 10005296
$ws.Range("K113").Value = 10005296
$ws.Range("M113").Value = -10003126

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 7055
$ws.Range("I34").Value = 416.66666
$ws.Range("J34").Value = 12033.75
$ws.Range("K34").Value = 1249.99998
$ws.Range("L34").Value = 36101.25
$ws.Range("M34").Value = -1165.99998
$ws.Range("N34").Value = -36269.25

$ws.Range("H44").Value = 4749.3
$ws.Range("I44").Value = 391.5
$ws.Range("J44").Value = 11286
$ws.Range("K44").Value = 1174.5
$ws.Range("L44").Value = 33858
$ws.Range("M44").Value = -776.5
$ws.Range("N44").Value = -34654

$ws.Range("H107").Value = 4337279.5
$ws.Range("J107").Value = 5356958.5
$ws.Range("L107").Value = 16070875.5
$ws.Range("N107").Value = -16074715.5

$ws.Range("H109").Value = 11267.77
$ws.Range("I109").Value = 450.2857
$ws.Range("K109").Value = 1350.8571
$ws.Range("M109").Value = -310.8571000000002

$ws.Range("H130").Value = 11397.2
$ws.Range("I130").Value = 5660
$ws.Range("K130").Value = 16980
$ws.Range("M130").Value = -11960

$ws.Range("H131").Value = 4269.737
$ws.Range("I131").Value = 2505.5386
$ws.Range("J131").Value = 8092.1665
$ws.Range("K131").Value = 7516.6158
$ws.Range("L131").Value = 24276.4995
$ws.Range("M131").Value = -2476.6158
$ws.Range("N131").Value = -34356.49950000001

$ws.Range("H139").Value = 2344.9033
$ws.Range("I139").Value = 1387.2778
$ws.Range("K139").Value = 4161.8334
$ws.Range("M139").Value = 978.1665999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2801.25
$ws.Range("I102").Value = 2792.4546
$ws.Range("K102").Value = 2792.4546
$ws.Range("M102").Value = -1170.4546

$ws.Range("H126").Value = 1766.0714
$ws.Range("I126").Value = 1742.0834
$ws.Range("J126").Value = 1910
$ws.Range("K126").Value = 5226.2502
$ws.Range("L126").Value = 5730
$ws.Range("M126").Value = -2756.2502
$ws.Range("N126").Value = -10670

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 11346.333
$ws.Range("I22").Value = 11849.833
$ws.Range("J22").Value = 9332.333000000001
$ws.Range("K22").Value = 11849.833
$ws.Range("L22").Value = 9332.333000000001
$ws.Range("M22").Value = -11554.833
$ws.Range("N22").Value = -9922.333000000001

$ws.Range("H27").Value = 11346.333
$ws.Range("I27").Value = 11849.833
$ws.Range("J27").Value = 9332.333000000001
$ws.Range("K27").Value = 11849.833
$ws.Range("L27").Value = 9332.333000000001
$ws.Range("M27").Value = -11742.833
$ws.Range("N27").Value = -9546.333000000001

$ws.Range("H34").Value = 54999.5
$ws.Range("I34").Value = 20000
$ws.Range("K34").Value = 20000
$ws.Range("M34").Value = -19828

$ws.Range("H55").Value = 1325.5385
$ws.Range("I55").Value = 1572.8334
$ws.Range("J55").Value = 1113.5714
$ws.Range("K55").Value = 1572.8334
$ws.Range("L55").Value = 1113.5714
$ws.Range("M55").Value = -1399.8334
$ws.Range("N55").Value = -1459.5714

$ws.Range("H122").Value = 3499.25
$ws.Range("I122").Value = 3357.465
$ws.Range("K122").Value = 10072.395
$ws.Range("M122").Value = -7622.395

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 61666.332
$ws.Range("I76").Value = 40000
$ws.Range("J76").Value = 72499.5
$ws.Range("K76").Value = 40000
$ws.Range("L76").Value = 72499.5
$ws.Range("M76").Value = -39685
$ws.Range("N76").Value = -73129.5

$ws.Range("H79").Value = 61666.332
$ws.Range("I79").Value = 40000
$ws.Range("J79").Value = 72499.5
$ws.Range("K79").Value = 40000
$ws.Range("L79").Value = 72499.5
$ws.Range("M79").Value = -38908
$ws.Range("N79").Value = -74683.5

$ws.Range("H113").Value = 653.2222
$ws.Range("I113").Value = 656
$ws.Range("K113").Value = 1968
$ws.Range("M113").Value = 202

$ws.Range("H117").Value = 96999.5
$ws.Range("J117").Value = 96999.5
$ws.Range("L117").Value = 96999.5
$ws.Range("N117").Value = -106177.5

$ws.Range("H132").Value = 402539.16
$ws.Range("I132").Value = 2710.647
$ws.Range("J132").Value = 1252174.8
$ws.Range("K132").Value = 8131.941
$ws.Range("L132").Value = 3756524.4
$ws.Range("M132").Value = -5601.941
$ws.Range("N132").Value = -3761584.4
